$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- sheet1 (LP1912) ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 13:18:40"
$ws1.Cells.Item(3,1).Value = "Total filas: 231"
$ws1.Cells.Item(90,1).Value = "07:24:45"
$ws1.Cells.Item(90,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(90,4).Value = 112
$ws1.Cells.Item(91,1).Value = "08:55:25"
$ws1.Cells.Item(91,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(91,4).Value = 21
$ws1.Cells.Item(120,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(121,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(133,1).Value = "09:26:30"
$ws1.Cells.Item(133,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(133,4).Value = 100
$ws1.Cells.Item(134,1).Value = "10:52:37"
$ws1.Cells.Item(134,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(134,4).Value = 14
$ws1.Cells.Item(162,1).Value = "10:52:37"
$ws1.Cells.Item(162,3).Value = "10_OLMOS"
$ws1.Cells.Item(162,4).Value = 74
$ws1.Cells.Item(163,1).Value = "10:13:53"
$ws1.Cells.Item(163,3).Value = "14_ABASTO"
$ws1.Cells.Item(163,4).Value = 113
$ws1.Cells.Item(164,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(178,1).Value = "11:17:39"
$ws1.Cells.Item(178,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(178,4).Value = 77
$ws1.Cells.Item(179,1).Value = "11:46:46"
$ws1.Cells.Item(179,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(179,4).Value = 48
$ws1.Cells.Item(193,1).Value = "12:50:41"
$ws1.Cells.Item(193,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(193,4).Value = 0
$ws1.Cells.Item(194,1).Value = "10:52:37"
$ws1.Cells.Item(194,3).Value = "15_ABASTO"
$ws1.Cells.Item(194,4).Value = 118
$ws1.Cells.Item(213,1).Value = "13:18:40"
$ws1.Cells.Item(213,2).Value = "13:36"
$ws1.Cells.Item(213,4).Value = 18
$ws1.Cells.Item(214,1).Value = "13:18:40"
$ws1.Cells.Item(214,2).Value = "13:36"
$ws1.Cells.Item(214,3).Value = "15_ABASTO"
$ws1.Cells.Item(214,4).Value = 18
$ws1.Cells.Item(215,1).Value = "12:35:30"
$ws1.Cells.Item(215,2).Value = "13:42"
$ws1.Cells.Item(215,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(215,4).Value = 67
$ws1.Cells.Item(216,1).Value = "13:18:40"
$ws1.Cells.Item(216,2).Value = "13:46"
$ws1.Cells.Item(216,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(216,4).Value = 28
$ws1.Cells.Item(217,2).Value = "13:46"
$ws1.Cells.Item(217,3).Value = "17_ROMERO"
$ws1.Cells.Item(217,4).Value = 105
$ws1.Cells.Item(218,1).Value = "12:50:41"
$ws1.Cells.Item(218,2).Value = "13:50"
$ws1.Cells.Item(218,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(218,4).Value = 60
$ws1.Cells.Item(219,2).Value = "13:50"
$ws1.Cells.Item(219,3).Value = "215A_EL PATO"
$ws1.Cells.Item(219,4).Value = 75
$ws1.Cells.Item(220,2).Value = "13:51"
$ws1.Cells.Item(220,3).Value = "215A_EL PATO"
$ws1.Cells.Item(220,4).Value = 110
$ws1.Cells.Item(221,1).Value = "12:01:11"
$ws1.Cells.Item(221,2).Value = "13:56"
$ws1.Cells.Item(221,3).Value = "225_GOMEZ"
$ws1.Cells.Item(221,4).Value = 115
$ws1.Cells.Item(222,1).Value = "12:35:30"
$ws1.Cells.Item(222,2).Value = "13:56"
$ws1.Cells.Item(222,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(222,4).Value = 81
$ws1.Cells.Item(223,1).Value = "12:01:11"
$ws1.Cells.Item(223,2).Value = "13:57"
$ws1.Cells.Item(223,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(223,4).Value = 116
$ws1.Cells.Item(224,1).Value = "13:18:40"
$ws1.Cells.Item(224,2).Value = "14:04"
$ws1.Cells.Item(224,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(224,4).Value = 46
$ws1.Cells.Item(225,2).Value = "14:04"
$ws1.Cells.Item(225,3).Value = "17_ROMERO"
$ws1.Cells.Item(225,4).Value = 89
$ws1.Cells.Item(226,1).Value = "12:50:41"
$ws1.Cells.Item(226,2).Value = "14:06"
$ws1.Cells.Item(226,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(226,4).Value = 76
$ws1.Cells.Item(227,1).Value = "12:35:30"
$ws1.Cells.Item(227,2).Value = "14:16"
$ws1.Cells.Item(227,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(227,4).Value = 101
$ws1.Cells.Item(228,1).Value = "12:50:41"
$ws1.Cells.Item(228,2).Value = "14:19"
$ws1.Cells.Item(228,3).Value = "215C_EL PATO"
$ws1.Cells.Item(228,4).Value = 89
$ws1.Cells.Item(228,5).Value = "LP1912"
$ws1.Cells.Item(229,1).Value = "12:35:30"
$ws1.Cells.Item(229,2).Value = "14:20"
$ws1.Cells.Item(229,3).Value = "215C_EL PATO"
$ws1.Cells.Item(229,4).Value = 105
$ws1.Cells.Item(229,5).Value = "LP1912"
$ws1.Cells.Item(230,1).Value = "12:35:30"
$ws1.Cells.Item(230,2).Value = "14:21"
$ws1.Cells.Item(230,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(230,4).Value = 106
$ws1.Cells.Item(230,5).Value = "LP1912"
$ws1.Cells.Item(231,1).Value = "12:50:41"
$ws1.Cells.Item(231,2).Value = "14:44"
$ws1.Cells.Item(231,3).Value = "14_ABASTO"
$ws1.Cells.Item(231,4).Value = 114
$ws1.Cells.Item(231,5).Value = "LP1912"
$ws1.Cells.Item(232,1).Value = "13:18:40"
$ws1.Cells.Item(232,2).Value = "14:56"
$ws1.Cells.Item(232,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(232,4).Value = 98
$ws1.Cells.Item(232,5).Value = "LP1912"
$ws1.Cells.Item(233,1).Value = "13:18:40"
$ws1.Cells.Item(233,2).Value = "14:58"
$ws1.Cells.Item(233,3).Value = "215B_EL PATO"
$ws1.Cells.Item(233,4).Value = 100
$ws1.Cells.Item(233,5).Value = "LP1912"
$ws1.Cells.Item(234,1).Value = "13:18:40"
$ws1.Cells.Item(234,2).Value = "15:00"
$ws1.Cells.Item(234,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(234,4).Value = 102
$ws1.Cells.Item(234,5).Value = "LP1912"
$ws1.Cells.Item(235,1).Value = "13:18:40"
$ws1.Cells.Item(235,2).Value = "15:05"
$ws1.Cells.Item(235,3).Value = "10_OLMOS"
$ws1.Cells.Item(235,4).Value = 107
$ws1.Cells.Item(235,5).Value = "LP1912"
$ws1.Cells.Item(236,1).Value = "13:18:40"
$ws1.Cells.Item(236,2).Value = "15:13"
$ws1.Cells.Item(236,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(236,4).Value = 115
$ws1.Cells.Item(236,5).Value = "LP1912"

# --- sheet2 (LP1912-215) ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 13:18:40"
$ws2.Cells.Item(3,1).Value = "Total filas: 29"
$ws2.Cells.Item(34,1).Value = "13:18:40"
$ws2.Cells.Item(34,2).Value = "14:58"
$ws2.Cells.Item(34,3).Value = "215B_EL PATO"
$ws2.Cells.Item(34,4).Value = 100
$ws2.Cells.Item(34,5).Value = "LP1912"

# --- sheet3 (6203-6173) ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 13:18:40"
$ws3.Cells.Item(3,1).Value = "Total filas: 35"
$ws3.Cells.Item(40,1).Value = "13:18:40"
$ws3.Cells.Item(40,2).Value = "14:52"
$ws3.Cells.Item(40,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(40,4).Value = 94
$ws3.Cells.Item(40,5).Value = "L6203"
